# Update cryptos list values (Price and Volume(1h)) per the Oct 7 2023 GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.910.04"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "1.632.69"
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'211.75"
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "'23.24"
$ws.Range("E8").Value = "  -0.94%  "
$ws.Range("E9").Value = "  -3.03%  "
$ws.Range("D10").Value = "'0.0613"
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("D11").Value = "'0.0880"
$ws.Range("E11").Value = "  +0.96%  "
$ws.Range("D12").Value = "1.863.85"
$ws.Range("E12").Value = "  -0.95%  "
$ws.Range("D13").Value = "1.637.01"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").Value = "'0.565"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").Value = "'65.26"
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").Value = "27.909.81"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").Value = "'230.19"
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("D19").Value = "0.0₃0722"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("E20").Value = "  -2.60%  "
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("D23").Value = "'10.37"
$ws.Range("E23").Value = "  -2.80%  "
$ws.Range("E24").Value = "  -4.27%  "
$ws.Range("D25").Value = "'154.31"
$ws.Range("E25").Value = "  +1.12%  "
$ws.Range("D26").Value = "'6.96"
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("D28").Value = "'15.62"
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  -0.98%  "
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("E33").Value = "  -0.69%  "
$ws.Range("D34").Value = "1.399.06"
$ws.Range("E34").Value = "  -3.42%  "
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("E36").Value = "  +9.52%  "
$ws.Range("E37").Value = "  +1.40%  "
$ws.Range("E38").Value = "  +0.37%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").Value = "'0.871"
$ws.Range("E40").Value = "  -2.14%  "
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").Value = "'66.82"
$ws.Range("E43").Value = "  -3.84%  "
$ws.Range("D44").Value = "'5.54"
$ws.Range("E44").Value = "  +2.48%  "
$ws.Range("E45").Value = "  +1.12%  "
$ws.Range("E46").Value = "  -1.12%  "
$ws.Range("D47").Value = "1.773.45"
$ws.Range("E47").Value = "  -0.99%  "
$ws.Range("D48").Value = "'87.65"
$ws.Range("E48").Value = "  -1.49%  "
$ws.Range("E49").Value = "  +0.58%  "
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("E51").Value = "  -0.27%  "
